# "Actualizaciones al modulo de reportes"
# Clears the dynamically-generated date/time stamp ("Fecha: ..." / "Hora: ...")
# out of the report header on Hoja1, and moves the active selection up to F8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# F9 held "Fecha: 20/03/2025" and F10 held "Hora: 11:28:15 Hrs." - strip the
# text but keep the cell formatting (style) intact.
$ws.Range("F9").ClearContents()
$ws.Range("F10").ClearContents()

# Reflect the new active cell/selection left after the edit.
$ws.Range("F8").Select()
